$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-row Price (D) and Volume/1h (E) updates ---
$ws.Range("D2").Value = "68.276.08"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "3.730.29"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'592.92"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'167.12"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "3.732.95"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").Value = "'6.46"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -6.37%  "
$ws.Range("D14").Value = "'36.22"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "4.355.71"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.727.75"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "68.174.93"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "'17.89"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("D19").Value = "'7.01"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").Value = "'466.76"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "'84.00"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "3.875.36"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'2.78"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").Value = "'29.86"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "'9.26"
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("D37").Value = "3.684.53"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  -11.41%  "
$ws.Range("D40").Value = "'0.137"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'0.995"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "'5.80"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "'1.94"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'43.24"
$ws.Range("E47").Value = "  +10.80%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  -0.15%  "

# --- Rows 50 & 51: Bittensor/Monero swap ranking positions (coin, link, price, volume all move) ---
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'146.25"
$ws.Range("E50").Value = "  +5.66%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'393.30"
$ws.Range("E51").Value = "  -0.83%  "
